# Add two new error-code rows for the merged phone/electronic address
# "primary contact updated" messages, mirroring the existing rows for
# the other address record-lifecycle messages on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A51").Value = 10049
$ws.Range("B51").Value = "message_10049_phone_address_primary_contact_updated"
$ws.Range("D51").Value = "Success"

$ws.Range("A52").Value = 10050
$ws.Range("B52").Value = "message_10050_electronic_address_primary_contact_updated"
$ws.Range("D52").Value = "Success"

# Restore the view to match the scrolled/selected state recorded after
# the edit (top-left cell + active selection).
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("B76").Select()
